$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh: update Price (D) and Volume(1h) (E) columns.
# Some new Price values parse as plain numbers; a leading apostrophe forces
# Excel to keep storing them as text (matching the original text cells)
# instead of converting them to numbers, and the style is reset right after
# so no extra number format sticks to the cell.
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '27.971.86'
$ws.Range('E2').Value = '  -3.20%  '
$ws.Range('D3').Value = '1.864.50'
$ws.Range('E3').Value = '  -2.18%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue $ws.Range('D5') '318.26'
$ws.Range('E5').Value = '  -1.82%  '
$ws.Range('E6').Value = '  +0.04%  '
Set-TextValue $ws.Range('D7') '0.4373'
$ws.Range('E7').Value = '  -4.65%  '
Set-TextValue $ws.Range('D8') '0.3698'
$ws.Range('E8').Value = '  -3.08%  '
Set-TextValue $ws.Range('D9') '0.07511'
$ws.Range('E9').Value = '  -2.51%  '
Set-TextValue $ws.Range('D10') '0.9366'
$ws.Range('E10').Value = '  -4.32%  '
Set-TextValue $ws.Range('D11') '21.26'
$ws.Range('D12').Value = '1.884.28'
$ws.Range('E12').Value = '  -0.12%  '
Set-TextValue $ws.Range('D13') '6.734'
$ws.Range('E13').Value = '  -3.06%  '
Set-TextValue $ws.Range('D14') '5.440'
$ws.Range('E14').Value = '  -4.15%  '
Set-TextValue $ws.Range('D15') '0.06831'
$ws.Range('E15').Value = '  -3.22%  '
$ws.Range('E16').Value = '  -0.04%  '
Set-TextValue $ws.Range('D17') '81.62'
$ws.Range('E17').Value = '  -2.54%  '
Set-TextValue $ws.Range('D18') '0.000009062'
$ws.Range('E18').Value = '  -4.07%  '
Set-TextValue $ws.Range('D19') '1.001'
$ws.Range('E19').Value = '  +0.06%  '
Set-TextValue $ws.Range('D20') '15.96'
$ws.Range('E20').Value = '  -4.04%  '
$ws.Range('D21').Value = '27.961.06'
$ws.Range('E21').Value = '  -3.21%  '
Set-TextValue $ws.Range('D22') '5.110'
$ws.Range('E22').Value = '  -3.70%  '
Set-TextValue $ws.Range('D23') '11.04'
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('D24').Value = '2.084.70'
$ws.Range('E24').Value = '  -1.37%  '
$ws.Range('E25').Value = '  -4.50%  '
Set-TextValue $ws.Range('D26') '154.07'
$ws.Range('E26').Value = '  -2.87%  '
$ws.Range('E28').Value = '  -3.83%  '
Set-TextValue $ws.Range('D29') '113.33'
$ws.Range('E29').Value = '  -3.43%  '
$ws.Range('E30').Value = '  -8.40%  '
Set-TextValue $ws.Range('D31') '0.08996'
$ws.Range('E31').Value = '  -3.32%  '
Set-TextValue $ws.Range('D32') '0.8106'
$ws.Range('E32').Value = '  -5.76%  '
Set-TextValue $ws.Range('D33') '4.816'
$ws.Range('E33').Value = '  -5.42%  '
$ws.Range('E34').Value = '  -5.46%  '
Set-TextValue $ws.Range('D35') '2.959'
$ws.Range('E35').Value = '  -2.68%  '
Set-TextValue $ws.Range('D36') '1.001'
$ws.Range('E36').Value = '  +0.01%  '
Set-TextValue $ws.Range('D37') '0.05495'
$ws.Range('E37').Value = '  -3.73%  '
$ws.Range('E38').Value = '  -3.64%  '
Set-TextValue $ws.Range('D39') '0.01979'
$ws.Range('E39').Value = '  -2.75%  '
Set-TextValue $ws.Range('D40') '2.921'
$ws.Range('E40').Value = '  -0.81%  '
Set-TextValue $ws.Range('D41') '0.5266'
$ws.Range('E41').Value = '  -4.07%  '
Set-TextValue $ws.Range('D42') '7.053'
$ws.Range('E42').Value = '  -5.36%  '
$ws.Range('E43').Value = '  -3.49%  '
Set-TextValue $ws.Range('D44') '8.796'
$ws.Range('E44').Value = '  -5.71%  '
Set-TextValue $ws.Range('D45') '0.06778'
$ws.Range('E45').Value = '  -1.66%  '
Set-TextValue $ws.Range('D46') '0.4907'
$ws.Range('E46').Value = '  -5.16%  '
Set-TextValue $ws.Range('D47') '10.63'
$ws.Range('E47').Value = '  -5.17%  '
Set-TextValue $ws.Range('D48') '106.55'
$ws.Range('E48').Value = '  -3.43%  '
$ws.Range('E49').Value = '  -5.37%  '
Set-TextValue $ws.Range('D50') '1.000'
$ws.Range('E50').Value = '  -0.07%  '
Set-TextValue $ws.Range('D51') '1.889'
$ws.Range('E51').Value = '  -12.20%  '
